$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I1").Value = 0.83889209981213608
$ws.Range("Q1").Value = 0.99718051199086699
$ws.Range("S1").Value = 0.82443854617857482
$ws.Range("R2").Value = 0.97234920188701945
$ws.Range("AM3").Value = 0.93064365753564071
$ws.Range("BM3").Value = 0.88582967468290696
$ws.Range("I4").Value = 0.86967017492176302
$ws.Range("N4").Value = 0.83498339382888909
$ws.Range("BO4").Value = 0.93318383974178254
$ws.Range("B5").Value = 0.79890565422218029
$ws.Range("BC5").Value = 0.99925400137425724
$ws.Range("D6").Value = 0.69235672719825458
$ws.Range("M6").Value = 0.80823070208948722
$ws.Range("E7").Value = 0.97897284390696981
$ws.Range("Y7").Value = 0.68521827402884083
$ws.Range("AM7").Value = 0.57944679258729437
$ws.Range("BM7").Value = 0.60792559920238198
$ws.Range("V8").Value = 0.8205877033210125
$ws.Range("AB8").Value = 0.8498102588915859
$ws.Range("AB9").Value = 0.85271961864728252
$ws.Range("AZ10").Value = 0.98557326951180924
$ws.Range("BA10").Value = 0.76969040978716652
$ws.Range("BI10").Value = 0.91601423857261821
$ws.Range("BG11").Value = 0.9698735785907705
$ws.Range("BJ11").Value = 0.97362362279438552
$ws.Range("Z12").Value = 0.91374084385486332
$ws.Range("BN12").Value = 0.95462961739345809
$ws.Range("AL13").Value = 0.71287489727275244
$ws.Range("BG13").Value = 0.77866662732573699
$ws.Range("BB14").Value = 0.77191394770368582
$ws.Range("N15").Value = 0.90061378155348115
$ws.Range("Z15").Value = 0.93360926146772716
$ws.Range("AF15").Value = 0.97680170396311672
$ws.Range("BB16").Value = 0.96493393313596565
$ws.Range("N17").Value = 0.97781304868050589
$ws.Range("O17").Value = 0.87023457809640836
$ws.Range("AL17").Value = 0.95450810024551647
$ws.Range("AZ18").Value = 0.86131589430878974
$ws.Range("B19").Value = 0.83511760804069501
$ws.Range("AD19").Value = 0.59696696588671816
$ws.Range("AA20").Value = 0.98852588520675111
$ws.Range("AX20").Value = 0.91402113351026937
$ws.Range("BF20").Value = 0.95506559283940851
$ws.Range("AY21").Value = 0.96237519393162407
$ws.Range("BC21").Value = 0.74791620179653617
$ws.Range("BE21").Value = 0.90515654148446278
$ws.Range("BI23").Value = 0.87442801761307787
$ws.Range("D24").Value = 0.9765173339765808
$ws.Range("L24").Value = 0.89465455156385443
$ws.Range("W24").Value = 0.85753508899413322
$ws.Range("Q25").Value = 0.9579611983698566
$ws.Range("H27").Value = 0.91557965833931765
$ws.Range("Y27").Value = 0.86745661959134002
$ws.Range("AV27").Value = 0.666982455024412
$ws.Range("C28").Value = 0.83130415354009402
$ws.Range("BD29").Value = 0.64266517380122146
$ws.Range("AC30").Value = 0.83013240890544671
$ws.Range("BJ30").Value = 0.88044020128143896
$ws.Range("AB31").Value = 0.92389359216454547
$ws.Range("AX31").Value = 0.94028616528026587
$ws.Range("S32").Value = 0.95984761509500438
$ws.Range("K33").Value = 0.94130326146522114
$ws.Range("D34").Value = 0.72097359220182733
$ws.Range("V34").Value = 0.87440907691248426
$ws.Range("AG34").Value = 0.92258256255944948
$ws.Range("AJ34").Value = 0.68705133399198193
$ws.Range("BP34").Value = 0.7697577873639736
$ws.Range("V35").Value = 0.82481491469238899
$ws.Range("AG35").Value = 0.99016134486733809
$ws.Range("AK35").Value = 0.87528178128691181
$ws.Range("AC36").Value = 0.96054328633323915
$ws.Range("AI36").Value = 0.93374725501777978
$ws.Range("AS37").Value = 0.97796538535743049
$ws.Range("BM37").Value = 0.86499530268398406
$ws.Range("AE38").Value = 0.98142952700351538
$ws.Range("H39").Value = 0.73803467882388518
$ws.Range("L40").Value = 0.7660196281467031
$ws.Range("AL40").Value = 0.96046480059082517
$ws.Range("AC41").Value = 0.91582084131791486
$ws.Range("AQ41").Value = 0.96106621187003749
$ws.Range("BN41").Value = 0.76798012693388584
$ws.Range("AG42").Value = 0.99239286138148852
$ws.Range("AX42").Value = 0.94184014441390929
$ws.Range("F43").Value = 0.91560319217482478
$ws.Range("H44").Value = 0.99782353174643657
$ws.Range("AM44").Value = 0.94370781846611318
$ws.Range("O45").Value = 0.88723684067540287
$ws.Range("BH45").Value = 0.73250406656285816
$ws.Range("BI45").Value = 0.83729092424259777
$ws.Range("AR46").Value = 0.98409832335600078
$ws.Range("AU46").Value = 0.9564855103028358
$ws.Range("BA46").Value = 0.88192613498248895
$ws.Range("BK46").Value = 0.9967323443777496
$ws.Range("M47").Value = 0.68194047589283158
$ws.Range("BG47").Value = 0.92751352806149123
$ws.Range("Q48").Value = 0.66816429522278487
$ws.Range("AT48").Value = 0.96621016597732146
$ws.Range("BF48").Value = 0.74570350993299706
$ws.Range("BF49").Value = 0.59794209605036752
$ws.Range("BG49").Value = 0.56208547384228003
$ws.Range("AT50").Value = 0.95873582194712392
$ws.Range("AZ51").Value = 0.91261283825726003
$ws.Range("Z52").Value = 0.87870818368467518
$ws.Range("BJ54").Value = 0.82737781798051602
$ws.Range("BM54").Value = 0.92528066003394982
$ws.Range("BN54").Value = 0.96172573028862818
$ws.Range("F55").Value = 0.78561556546906364
$ws.Range("AI55").Value = 0.71483503365991563
$ws.Range("AT55").Value = 0.93585007441670243
$ws.Range("BC56").Value = 0.81363903977487873
$ws.Range("BG56").Value = 0.74252280992909125
$ws.Range("BI56").Value = 0.82041348799337555
$ws.Range("AH57").Value = 0.80397525966343408
$ws.Range("BP57").Value = 0.93193685994814524
$ws.Range("I58").Value = 0.92583499457392326
$ws.Range("R58").Value = 0.98099936710720026
$ws.Range("Z58").Value = 0.78344486672628377
$ws.Range("BE58").Value = 0.78013577480177232
$ws.Range("BK58").Value = 0.88471521160976496
$ws.Range("AD59").Value = 0.96642864775296655
$ws.Range("F60").Value = 0.97574715103716914
$ws.Range("AC61").Value = 0.80409201898332605
$ws.Range("B62").Value = 0.98431425607162359
$ws.Range("O62").Value = 0.97391670310456102
$ws.Range("P63").Value = 0.90758485259106569
$ws.Range("AC63").Value = 0.82494846059754146
$ws.Range("BH64").Value = 0.76301328634830057
$ws.Range("AN65").Value = 0.82163811074945359
$ws.Range("J66").Value = 0.81053593307634653
$ws.Range("W66").Value = 0.98148306245409733
$ws.Range("BL66").Value = 0.87350468004802506
$ws.Range("BO66").Value = 0.83814111311153605
$ws.Range("AC67").Value = 0.83927350131880885
$ws.Range("AV67").Value = 0.99433214040815932
$ws.Range("BP67").Value = 0.94078331758237221
$ws.Range("F68").Value = 0.7213945096889125
